# Update mail-merge field results (cached values) in the document.
# Mapping of MERGEFIELD code name -> new result text.
$d = $word.ActiveDocument

$updates = @{
    "MERGEFIELD NO"     = "E30"
    "MERGEFIELD NAMA"   = "I MADE SURYA D."
    "MERGEFIELD SEPATU" = "46"
    "MERGEFIELD UBN_1"  = "53"
    "MERGEFIELD UBN_2"  = "63"
    "MERGEFIELD UBN_5"  = "116"
    "MERGEFIELD UBN_6"  = "112"
    "MERGEFIELD UBN_7"  = "120"
    "MERGEFIELD UBN_8"  = "75"
    "MERGEFIELD UBN_9"  = "45"
    "MERGEFIELD UH_1"   = "53"
    "MERGEFIELD UH_2"   = "64"
    "MERGEFIELD UH_5"   = "31"
    "MERGEFIELD UH_6"   = "29"
    "MERGEFIELD UH_7"   = "31"
    "MERGEFIELD UH_8"   = "74"
    "MERGEFIELD UH_9"   = "45"
}

foreach ($f in $d.Fields) {
    $code = $f.Code.Text.Trim()
    if ($updates.ContainsKey($code)) {
        $r = $f.Result
        $rng = $d.Range($r.Start, $r.End)
        $rng.Text = $updates[$code]
    }
}
